$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row 9: "Digi-Key Order" -> "Digi-Key Order 1" (label + linked pdf)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Digi-Key Order 1"
$ws.Range("D9").Value = "Digi-Key Order 1.pdf"

# ---------------------------------------------------------------------------
# 2. Clear the old "Total" row (13) and old footnote row (15); they move down.
# ---------------------------------------------------------------------------
$ws.Range("D13:E13").ClearContents()
$ws.Range("A15").ClearContents()

# ---------------------------------------------------------------------------
# 3. New receipt rows: Banggood Order 1 (12), Hobby King Order 3 (13),
#    Digi-Key Order 2 (14, still awaiting a receipt so no cost/link yet).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Banggood Order 1"
$ws.Range("B12").Value = 43092
$ws.Range("C12").Value = "Brian"
$ws.Range("D12").Value = "Banggood Order 1.pdf"
$ws.Range("E12").Value = 23.08
$ws.Range("F12").Value = "120W AC/DC Power Adapter"

$ws.Range("A13").Value = "Hobby King Order 3"
$ws.Range("B13").Value = 43096
$ws.Range("C13").Value = "Brian"
$ws.Range("D13").Value = "Hobby King Order 3.pdf"
$ws.Range("E13").Value = 74.79
$ws.Range("F13").Value = "Higher Voltage, Low KV Motors"

$ws.Range("A14").Value = "Digi-Key Order 2"
$ws.Range("B14").Value = 43097
$ws.Range("C14").Value = "Brian"
$ws.Range("D14").Style = "Hyperlink"
$ws.Range("F14").Value = "5.5x2.5mm Barrel Jacks"

# ---------------------------------------------------------------------------
# 4. New "Total" row (16) and footnote row (18).
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = "Total"
$ws.Range("D16").Font.Bold = $true
$ws.Range("E16").Formula = "=SUM(E2:E15)"

$ws.Range("A18").Value = "*Not sure if Feron will fund this order"

# ---------------------------------------------------------------------------
# 5. Rebuild the hyperlinks. The engine's Hyperlinks.Add doesn't replace an
#    existing hyperlink on a cell (it stacks a duplicate entry), so drop all
#    of them and recreate the full, correctly-ordered set in one pass.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = @(
  @("D2",  "Amazon%20Order%201.pdf"),
  @("D3",  "Hobby%20King%20Order%201.pdf"),
  @("D4",  "Amazon%20Order%202.pdf"),
  @("D5",  "Amazon%20Order%203.pdf"),
  @("D6",  "Amazon%20Order%204.pdf"),
  @("D7",  "Hobby%20King%20Order%202.pdf"),
  @("D8",  "Amazon%20Order%205.pdf"),
  @("D9",  "Digi-Key%20Order%201.pdf"),
  @("D10", "Arrow%20Order.pdf"),
  @("D11", "Amazon%20Order%206.pdf"),
  @("D12", "Banggood%20Order%201.pdf"),
  @("D13", "Hobby%20King%20Order%203.pdf")
)

foreach ($pair in $links) {
    $cellref = $pair[0]
    $target = $pair[1]
    $range = $ws.Range($cellref)
    $ws.Hyperlinks.Add($range, $target, "", "", $range.Value())
}

# Re-adding the hyperlinks above resets each cell to a freshly minted style;
# put them back on the shared "Hyperlink" cell style used throughout the file.
foreach ($pair in $links) {
    $ws.Range($pair[0]).Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 6. Selection, matching the workbook as last saved by the author.
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
